$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (door / location annotations) ---
# Order matters: it controls the order new entries are appended to the
# shared-strings table, so we touch the cells in the same order the
# target workbook's sharedStrings.xml lists its new unique strings.

$ws.Range("C10").Value = "door3 / (18,0,4)"

$ws.Range("C12").Value = "door5(28,0,4)"

$ws.Range("C11").Value = "door6 (38,0,7)"
# Colour the "(38,0,7)" portion of C11 red, leaving "door6 " in the default colour.
$chars = $ws.Range("C11").Characters(7, 8)
$chars.Font.Color = 255

$ws.Range("C14").Value = "door4(40,0,6)"

$ws.Range("C15").Value = "door6(34,0,8)"

$ws.Range("C16").Value = "door6(37,0,7)"

$ws.Range("C17").Value = "door3(7,0,15)"

$ws.Range("C18").Value = "door6(34,0,8)"

$ws.Range("C20").Value = "doorEntrance (43,00,14)"

$ws.Range("A21").Value = "second level"
$ws.Range("C21").Value = "doorKey4(67,0,77)"

# --- Page setup (Page Setup dialog: paper size + orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- View state: zoom + scroll position + active selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 94
$ws.Range("D17").Select()
